$d = $word.ActiveDocument

# 1) Merge "Collaboration" + "/" runs into "Collaboration/" (no visible text
#    change, but normalize the adjacent "Collaboration" "/" run split so the
#    resulting OOXML has a single run for "Collaboration/").
$d.Content.Find.Execute("Collaboration/Teamwork", $true, $false, $false, $false, $false, $true, 1, $false, "Collaboration/Teamwork", 2) | Out-Null

# 2) Extend the harassment clause: drop the trailing period and append the
#    new consequence clause.
$d.Content.Find.Execute("and any other personal characteristics.", $true, $false, $false, $false, $false, $true, 1, $false, "and any other personal characteristics otherwise lead to severe consequences", 2) | Out-Null

# 3) Replace "abide to them" with "follow them".
$d.Content.Find.Execute("abide to them", $true, $false, $false, $false, $false, $true, 1, $false, "follow them", 2) | Out-Null
